# Update "想去人数" (interest count) values in the F column across the
# workbook's four sheets, matching the data refresh captured in the diff.

$wb = $excel.ActiveWorkbook

$sheetExhibition = $wb.Worksheets.Item("展览")
$sheetShow       = $wb.Worksheets.Item("演出")
$sheetLocalLife  = $wb.Worksheets.Item("本地生活")
$sheetAll        = $wb.Worksheets.Item("全部类型")

# 展览 (Exhibition) sheet
$sheetExhibition.Range("F6").Value  = 362
$sheetExhibition.Range("F8").Value  = 146
$sheetExhibition.Range("F10").Value = 214
$sheetExhibition.Range("F11").Value = 5912
$sheetExhibition.Range("F13").Value = 41
$sheetExhibition.Range("F14").Value = 491
$sheetExhibition.Range("F17").Value = 355
$sheetExhibition.Range("F22").Value = 131
$sheetExhibition.Range("F25").Value = 1017
$sheetExhibition.Range("F27").Value = 1797
$sheetExhibition.Range("F28").Value = 462

# 演出 (Performance) sheet
$sheetShow.Range("F5").Value = 267
$sheetShow.Range("F8").Value = 47

# 本地生活 (Local Life) sheet
$sheetLocalLife.Range("F2").Value = 221

# 全部类型 (All Types) sheet - aggregated view of the other sheets
$sheetAll.Range("F2").Value  = 221
$sheetAll.Range("F8").Value  = 362
$sheetAll.Range("F10").Value = 146
$sheetAll.Range("F12").Value = 214
$sheetAll.Range("F13").Value = 5912
$sheetAll.Range("F15").Value = 41
$sheetAll.Range("F17").Value = 491
$sheetAll.Range("F20").Value = 355
$sheetAll.Range("F25").Value = 267
$sheetAll.Range("F29").Value = 47
$sheetAll.Range("F32").Value = 131
$sheetAll.Range("F35").Value = 1017
$sheetAll.Range("F37").Value = 1797
$sheetAll.Range("F38").Value = 462
